# Climate Change Notes.docx
#
# Fixes the typo "do not  have nerves, they can't move" (double space
# between "not" and "have") -> "do not have nerves, they can't move".
#
# The visible text only loses one redundant space, but because the
# original sentence was already split into many small <w:r> runs (left
# over from earlier revisions, with a <w:proofErr> gramStart/gramEnd
# pair bracketing "not  have"), re-typing the area re-buckets those
# runs: "not have"/" nerves"/", they " become fresh runs in front of the
# gramStart marker, while the gramStart/gramEnd pair ends up wrapping a
# lone "can't" run, and the remaining tail becomes " move".
#
# Strategy:
#   1. Replace the text physically bracketed by gramStart/gramEnd
#      ("not  have") with "can't" - this keeps the proofErr markers
#      anchored to that run.
#   2. Insert "not have nerves, they " immediately before that run (so
#      it ends up before gramStart, after "do ").
#   3. Delete the now-duplicated " nerves, they can't" that got left
#      behind right after gramEnd, leaving " move" in place.
#   4. Re-split the runs that the engine coalesced while doing the
#      above back into the pieces the final document uses, by toggling
#      Bold on/off (a formatting no-op) across each desired run's span -
#      this forces a run boundary without altering the text.

$d = $word.ActiveDocument

# --- Step 1: turn the gramStart/gramEnd-wrapped "not  have" into "can't" ---
$d.Content.Find.Execute("not  have", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "can’t", 2) | Out-Null

# --- Step 2: insert the rest of the corrected phrase right before it ---
$text = $d.Content.Text
$idx = $text.IndexOf("can’t")
$insertionPoint = $d.Range($idx, $idx)
$insertionPoint.InsertBefore("not have nerves, they ")

# --- Step 3: drop the stale tail copy that used to sit after gramEnd ---
$text = $d.Content.Text
$idx = $text.IndexOf("can’t")
$tailStart = $idx + ([string]"can’t").Length
$stale = " nerves, they can’t"
$staleRange = $d.Range($tailStart, $tailStart + $stale.Length)
$staleRange.Text = ""

# --- Step 4: re-split the coalesced runs to match the target layout ---
function Split-RunSpans($startPos, [object[]]$spans) {
    foreach ($span in $spans) {
        $s = $startPos + $span[0]
        $e = $startPos + $span[1]
        $r = $d.Range($s, $e)
        $r.Bold = 1
        $r.Bold = 0
    }
}

# "plants do not have nerves, they " -> "plants "/"do "/"not have"/" nerves"/", they "
$text = $d.Content.Text
$run1Start = $text.IndexOf("plants do not have nerves, they ")
Split-RunSpans $run1Start @(
    , @(0, 7)   # "plants "
    , @(7, 10)  # "do "
    , @(10, 18) # "not have"
    , @(18, 25) # " nerves"
    , @(25, 32) # ", they "
)

# " move meaning that if climate change continues, plants will die." ->
#   " move" / " meaning that if climate change continues" / ", plants will die."
$text = $d.Content.Text
$run2Start = $text.IndexOf(" move meaning that if climate change continues, plants will die.")
Split-RunSpans $run2Start @(
    , @(0, 5)    # " move"
    , @(5, 46)   # " meaning that if climate change continues"
    , @(46, 65)  # ", plants will die."
)

Write-Output "done"
